$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.6
$ws.Range("Q2").Value = 1.89
$ws.Range("R2").Value = 2.01
$ws.Range("AI2").Value = 29

# Row 3
$ws.Range("G3").Value = 1.73
$ws.Range("H3").Value = 3.25
$ws.Range("I3").Value = 5.75
$ws.Range("J3").Value = 2.5
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 5
$ws.Range("AI3").Value = 26
$ws.Range("AO3").Value = 10
$ws.Range("AP3").Value = 29

# Row 4
$ws.Range("G4").Value = 1.8
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.5
$ws.Range("L4").Value = 5
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.25
$ws.Range("R4").Value = 1.62
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1.73
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 7.5
$ws.Range("Z4").Value = 13
$ws.Range("AC4").Value = 7.5
$ws.Range("AE4").Value = 19
$ws.Range("AF4").Value = 67
$ws.Range("AH4").Value = 11
$ws.Range("AI4").Value = 23
$ws.Range("AJ4").Value = 17
$ws.Range("AK4").Value = 51
$ws.Range("AL4").Value = 41
$ws.Range("AN4").Value = 3.6
$ws.Range("AO4").Value = 9.5
$ws.Range("AU4").Value = 9
$ws.Range("AV4").Value = 67
$ws.Range("AW4").Value = 6.5
$ws.Range("AX4").Value = 26
$ws.Range("AZ4").Value = 101
$ws.Range("BA4").Value = 126

# Row 10
$ws.Range("Q10").Value = 2.25
$ws.Range("R10").Value = 1.62
$ws.Range("BD10").Value = 151
